$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Bharat"
$ws.Range("B1").Value = "Absent"

$ws.Range("C1").Value = 45480
$ws.Range("C1").NumberFormat = "yyyy-mm-dd"
